$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.094.50"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "3.816.22"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'701.20"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").Value = "'172.20"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "3.814.93"
$ws.Range("E7").Value = "  -0.75%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("E10").Value = "  -0.63%  "

$ws.Range("E11").Value = "  +2.70%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "4.459.68"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "3.840.65"
$ws.Range("E16").Value = "  +3.89%  "

$ws.Range("D17").Value = "71.149.54"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "'509.73"
$ws.Range("E21").Value = "  +3.46%  "

$ws.Range("D22").Value = "'10.72"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("E23").Value = "  +0.53%  "

$ws.Range("D24").Value = "'84.09"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").Value = "3.969.32"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "'12.05"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -4.34%  "

$ws.Range("D31").Value = "'3.02"
$ws.Range("E31").Value = "  -4.40%  "

$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "'2.24"
$ws.Range("E33").Value = "  -1.16%  "

$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("E35").Value = "  -5.01%  "

$ws.Range("E36").Value = "  +0.47%  "

$ws.Range("D37").Value = "3.778.65"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("D40").Value = "'2.37"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("D43").Value = "'3.29"
$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").Value = "'171.08"
$ws.Range("E45").Value = "  +4.70%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").Value = "'0.000311"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("D48").Value = "'49.45"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("D49").Value = "'427.22"
$ws.Range("E49").Value = "  +4.97%  "

$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.294"
$ws.Range("E51").Value = "  -1.52%  "
